# Delete the second data row (worksheet row 2), which contained the
# "com.singleton.strechy / taxi game / shamirnaftali@gmail.com / irisalmog47@gmail.com"
# review entry. All rows below it shift up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Delete()

# Update the active selection to match the post-edit state (B2).
$ws.Range("B2").Select()
